$d = $word.ActiveDocument

# wdReplace constants
$wdReplaceNone = 0
$wdReplaceOne  = 1
$wdReplaceAll  = 2
$wdFindContinue = 1

# ------------------------------------------------------------------
# 1) Title: "Play Lucky U Free: Review of Mechanics & RTP Rate"
#    -> "Play Lucky U for Free"
#    (occurs twice: the H1 heading and the bold run near the end;
#     ReplaceAll handles both in one pass)
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Play Lucky U Free: Review of Mechanics & RTP Rate",
    $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "Play Lucky U for Free", $wdReplaceAll) | Out-Null

# ------------------------------------------------------------------
# 2) "What we like" list: insert a new bullet "Easy to understand
#    gameplay mechanics" right before "Impressive 96% RTP rate"
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Impressive 96% RTP rate*") {
        $p.Range.InsertParagraphBefore()
        $newPara = $d.Paragraphs.Item($i)
        $newPara.Range.Text = "Easy to understand gameplay mechanics"
        break
    }
}

# ------------------------------------------------------------------
# 3) "Lucky You mode with free spins and multipliers"
#    -> "Lucky You mode with free spins and multiplied wins"
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Lucky You mode with free spins and multipliers",
    $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "Lucky You mode with free spins and multiplied wins", $wdReplaceAll) | Out-Null

# ------------------------------------------------------------------
# 4) Remove the "Engaging gameplay mechanics" bullet entirely
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Engaging gameplay mechanics*") {
        $p.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 5) "Aesthetic appeal with unique symbols"
#    -> "Visually appealing symbol design"
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Aesthetic appeal with unique symbols",
    $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "Visually appealing symbol design", $wdReplaceAll) | Out-Null

# ------------------------------------------------------------------
# 6) "What we don't like" list
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Limited duration of the Lucky You mode",
    $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "Lucky You mode is short-lived", $wdReplaceAll) | Out-Null

$d.Content.Find.Execute(
    "No progressive jackpot",
    $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "Limited variety in symbols", $wdReplaceAll) | Out-Null

# ------------------------------------------------------------------
# 7) Meta description (italic run)
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "Read our review of Lucky U, the Playtech online slot game featuring impressive RTP, Lucky You mode, and appealing symbol design. Play for free today.",
    $true, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "Read our review of Lucky U and play this engaging slot game for free.", $wdReplaceAll) | Out-Null

Write-Output "Done"
